# "preparing CPU for testing" — wired CPU / test ROM prep
#
# Major oversight found in the AL (arithmetic/logic) operations: the
# destination-register field documented in the "Extra info" column used
# bits a[15..12] for every AL op, but the actual instruction encoding
# (see column C / B) reserves a[15..12] for selecting the FIRST SOURCE
# operand register, while the destination is really chosen by a[11..8].
# Fix up the documentation text for each affected instruction, and correct
# the str() signature (operand order is reg/imm, not imm/reg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# str(reg/imm, address in mem) - operand order correction
$ws.Range("A6").Value = "str(reg/imm, address in mem)   "

# AL instructions: destination register is a[11..8], not a[15..12]
$ws.Range("D14").Value = "reg( a[11..8] ) <- a[11..0] / reg( a[15..12] ) + b[11..0] / reg( b[15..12] )"
$ws.Range("D15").Value = "reg( a[11..8] ) <- a[11..0] / reg( a[15..12] ) - b[11..0] / reg( b[15..12] )"
$ws.Range("D16").Value = "reg( a[11..8] ) <- a[11..0] / reg( a[15..12] ) * b[11..0] / reg( b[15..12] )"
$ws.Range("D17").Value = "reg( a[11..8] ) <- reg( a[15..12] ) >> b"
$ws.Range("D18").Value = "reg( a[11..8] ) <- a[11..0] / reg( a[15..12] ) & b[11..0] / reg( b[15..12] )"
$ws.Range("D19").Value = "reg( a[11..8] ) <- a[11..0] / reg( a[15..12] ) | b[11..0] / reg( b[15..12] )"
$ws.Range("D20").Value = "reg( a[11..8] ) <- ! b[11..0]"

# move the live selection to where work left off while writing test ROM
$null = $ws.Range("D23").Select()
